# This script updates the "dSF" column (F) values for specific rows to
# reflect a repull/recalculation of data, per the commit message
# "repull data, push all data, mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column F
$updates = @{
    5  = 0
    9  = 1
    11 = 1
    13 = -2
    18 = 2
    19 = 2
    26 = -11
    30 = 0
    31 = -2
    42 = -2
    46 = -1
    48 = -1
    56 = 3
    58 = 6
    59 = -4
    63 = -3
    64 = -3
    65 = -1
    67 = -7
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
